# Update "想去人数" (F) / "最低票价" (G) figures for 苏州-漫展信息.xlsx
# on both the "展览" and "全部类型" sheets, matching the latest scrape.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5213
$ws1.Range("F6").Value = 206
$ws1.Range("G7").Value = 50
$ws1.Range("F11").Value = 68
$ws1.Range("F13").Value = 37
$ws1.Range("F14").Value = 6602
$ws1.Range("F15").Value = 36
$ws1.Range("F20").Value = 1043
$ws1.Range("F21").Value = 15815
$ws1.Range("F22").Value = 1561
$ws1.Range("F23").Value = 25
$ws1.Range("F24").Value = 309
$ws1.Range("F27").Value = 11192
$ws1.Range("F28").Value = 806
$ws1.Range("F30").Value = 275
$ws1.Range("F31").Value = 380
$ws1.Range("F32").Value = 29

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5213
$ws4.Range("F6").Value = 206
$ws4.Range("G7").Value = 50
$ws4.Range("F12").Value = 68
$ws4.Range("F15").Value = 37
$ws4.Range("F16").Value = 6602
$ws4.Range("F17").Value = 36
$ws4.Range("F23").Value = 1043
$ws4.Range("F24").Value = 15815
$ws4.Range("F25").Value = 1561
$ws4.Range("F26").Value = 25
$ws4.Range("F27").Value = 309
$ws4.Range("F31").Value = 11192
$ws4.Range("F32").Value = 806
$ws4.Range("F34").Value = 275
$ws4.Range("F35").Value = 380
$ws4.Range("F36").Value = 29
